$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scattered single-cell updates (column C only) ---
$ws.Range("C22").Value = 3
$ws.Range("C31").Value = 3
$ws.Range("C1008").Value = 46
$ws.Range("C1023").Value = 56
$ws.Range("C1092").Value = 48
$ws.Range("C1126").Value = 57
$ws.Range("C1131").Value = 69
$ws.Range("C1219").Value = 32
$ws.Range("C1236").Value = 37
$ws.Range("C1246").Value = 37
$ws.Range("C1251").Value = 41
$ws.Range("C1266").Value = 29
$ws.Range("C1284").Value = 41
$ws.Range("C1296").Value = 13
$ws.Range("C1303").Value = 31
$ws.Range("C1307").Value = 4
$ws.Range("C1309").Value = 26

# --- Bulk refresh of rows 1314-1408 (data revision + 5 new appended rows) ---
$data = New-Object 'object[,]' 95,3
$data[0,0] = 44219
$data[0,1] = "50-59"
$data[0,2] = 4
$data[1,0] = 44219
$data[1,1] = "60-69"
$data[1,2] = 10
$data[2,0] = 44219
$data[2,1] = "70-79"
$data[2,2] = 12
$data[3,0] = 44219
$data[3,1] = "80+"
$data[3,2] = 25
$data[4,0] = 44220
$data[4,1] = "50-59"
$data[4,2] = 5
$data[5,0] = 44220
$data[5,1] = "60-69"
$data[5,2] = 7
$data[6,0] = 44220
$data[6,1] = "70-79"
$data[6,2] = 12
$data[7,0] = 44220
$data[7,1] = "80+"
$data[7,2] = 26
$data[8,0] = 44221
$data[8,1] = "40-49"
$data[8,2] = 1
$data[9,0] = 44221
$data[9,1] = "50-59"
$data[9,2] = 1
$data[10,0] = 44221
$data[10,1] = "60-69"
$data[10,2] = 11
$data[11,0] = 44221
$data[11,1] = "70-79"
$data[11,2] = 12
$data[12,0] = 44221
$data[12,1] = "80+"
$data[12,2] = 21
$data[13,0] = 44222
$data[13,1] = "0-19"
$data[13,2] = 1
$data[14,0] = 44222
$data[14,1] = "50-59"
$data[14,2] = 7
$data[15,0] = 44222
$data[15,1] = "60-69"
$data[15,2] = 7
$data[16,0] = 44222
$data[16,1] = "70-79"
$data[16,2] = 11
$data[17,0] = 44222
$data[17,1] = "80+"
$data[17,2] = 26
$data[18,0] = 44223
$data[18,1] = "40-49"
$data[18,2] = 1
$data[19,0] = 44223
$data[19,1] = "50-59"
$data[19,2] = 3
$data[20,0] = 44223
$data[20,1] = "60-69"
$data[20,2] = 9
$data[21,0] = 44223
$data[21,1] = "70-79"
$data[21,2] = 18
$data[22,0] = 44223
$data[22,1] = "80+"
$data[22,2] = 24
$data[23,0] = 44224
$data[23,1] = "60-69"
$data[23,2] = 7
$data[24,0] = 44224
$data[24,1] = "70-79"
$data[24,2] = 16
$data[25,0] = 44224
$data[25,1] = "80+"
$data[25,2] = 24
$data[26,0] = 44225
$data[26,1] = "30-39"
$data[26,2] = 1
$data[27,0] = 44225
$data[27,1] = "50-59"
$data[27,2] = 4
$data[28,0] = 44225
$data[28,1] = "60-69"
$data[28,2] = 6
$data[29,0] = 44225
$data[29,1] = "70-79"
$data[29,2] = 11
$data[30,0] = 44225
$data[30,1] = "80+"
$data[30,2] = 22
$data[31,0] = 44226
$data[31,1] = "20-29"
$data[31,2] = 1
$data[32,0] = 44226
$data[32,1] = "50-59"
$data[32,2] = 5
$data[33,0] = 44226
$data[33,1] = "60-69"
$data[33,2] = 4
$data[34,0] = 44226
$data[34,1] = "70-79"
$data[34,2] = 19
$data[35,0] = 44226
$data[35,1] = "80+"
$data[35,2] = 17
$data[36,0] = 44227
$data[36,1] = "40-49"
$data[36,2] = 3
$data[37,0] = 44227
$data[37,1] = "50-59"
$data[37,2] = 2
$data[38,0] = 44227
$data[38,1] = "60-69"
$data[38,2] = 11
$data[39,0] = 44227
$data[39,1] = "70-79"
$data[39,2] = 12
$data[40,0] = 44227
$data[40,1] = "80+"
$data[40,2] = 11
$data[41,0] = 44228
$data[41,1] = "20-29"
$data[41,2] = 1
$data[42,0] = 44228
$data[42,1] = "40-49"
$data[42,2] = 1
$data[43,0] = 44228
$data[43,1] = "50-59"
$data[43,2] = 3
$data[44,0] = 44228
$data[44,1] = "60-69"
$data[44,2] = 8
$data[45,0] = 44228
$data[45,1] = "70-79"
$data[45,2] = 14
$data[46,0] = 44228
$data[46,1] = "80+"
$data[46,2] = 18
$data[47,0] = 44229
$data[47,1] = "30-39"
$data[47,2] = 1
$data[48,0] = 44229
$data[48,1] = "40-49"
$data[48,2] = 3
$data[49,0] = 44229
$data[49,1] = "50-59"
$data[49,2] = 2
$data[50,0] = 44229
$data[50,1] = "60-69"
$data[50,2] = 5
$data[51,0] = 44229
$data[51,1] = "70-79"
$data[51,2] = 14
$data[52,0] = 44229
$data[52,1] = "80+"
$data[52,2] = 22
$data[53,0] = 44230
$data[53,1] = "50-59"
$data[53,2] = 3
$data[54,0] = 44230
$data[54,1] = "60-69"
$data[54,2] = 11
$data[55,0] = 44230
$data[55,1] = "70-79"
$data[55,2] = 10
$data[56,0] = 44230
$data[56,1] = "80+"
$data[56,2] = 25
$data[57,0] = 44231
$data[57,1] = "20-29"
$data[57,2] = 1
$data[58,0] = 44231
$data[58,1] = "50-59"
$data[58,2] = 3
$data[59,0] = 44231
$data[59,1] = "60-69"
$data[59,2] = 5
$data[60,0] = 44231
$data[60,1] = "70-79"
$data[60,2] = 6
$data[61,0] = 44231
$data[61,1] = "80+"
$data[61,2] = 19
$data[62,0] = 44232
$data[62,1] = "70-79"
$data[62,2] = 7
$data[63,0] = 44232
$data[63,1] = "80+"
$data[63,2] = 18
$data[64,0] = 44233
$data[64,1] = "40-49"
$data[64,2] = 1
$data[65,0] = 44233
$data[65,1] = "50-59"
$data[65,2] = 1
$data[66,0] = 44233
$data[66,1] = "60-69"
$data[66,2] = 6
$data[67,0] = 44233
$data[67,1] = "70-79"
$data[67,2] = 7
$data[68,0] = 44233
$data[68,1] = "80+"
$data[68,2] = 9
$data[69,0] = 44234
$data[69,1] = "50-59"
$data[69,2] = 2
$data[70,0] = 44234
$data[70,1] = "60-69"
$data[70,2] = 5
$data[71,0] = 44234
$data[71,1] = "70-79"
$data[71,2] = 6
$data[72,0] = 44234
$data[72,1] = "80+"
$data[72,2] = 13
$data[73,0] = 44235
$data[73,1] = "60-69"
$data[73,2] = 3
$data[74,0] = 44235
$data[74,1] = "70-79"
$data[74,2] = 4
$data[75,0] = 44235
$data[75,1] = "80+"
$data[75,2] = 11
$data[76,0] = 44236
$data[76,1] = "40-49"
$data[76,2] = 1
$data[77,0] = 44236
$data[77,1] = "50-59"
$data[77,2] = 5
$data[78,0] = 44236
$data[78,1] = "60-69"
$data[78,2] = 5
$data[79,0] = 44236
$data[79,1] = "70-79"
$data[79,2] = 3
$data[80,0] = 44236
$data[80,1] = "80+"
$data[80,2] = 9
$data[81,0] = 44237
$data[81,1] = "40-49"
$data[81,2] = 1
$data[82,0] = 44237
$data[82,1] = "50-59"
$data[82,2] = 1
$data[83,0] = 44237
$data[83,1] = "60-69"
$data[83,2] = 3
$data[84,0] = 44237
$data[84,1] = "70-79"
$data[84,2] = 10
$data[85,0] = 44237
$data[85,1] = "80+"
$data[85,2] = 3
$data[86,0] = 44238
$data[86,1] = "50-59"
$data[86,2] = 1
$data[87,0] = 44238
$data[87,1] = "60-69"
$data[87,2] = 4
$data[88,0] = 44238
$data[88,1] = "70-79"
$data[88,2] = 4
$data[89,0] = 44238
$data[89,1] = "80+"
$data[89,2] = 3
$data[90,0] = 44239
$data[90,1] = "60-69"
$data[90,2] = 1
$data[91,0] = 44239
$data[91,1] = "70-79"
$data[91,2] = 3
$data[92,0] = 44239
$data[92,1] = "80+"
$data[92,2] = 7
$data[93,0] = 44240
$data[93,1] = "60-69"
$data[93,2] = 1
$data[94,0] = 44240
$data[94,1] = "80+"
$data[94,2] = 2
$ws.Range("A1314:C1408").Value = $data

# Ensure new date cells (previously empty rows) carry the same date number format as the rest of column A
$ws.Range("A1404:A1408").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count()
